# Edit script: update B..F and I..N value columns in rows 2-25 of Sheet1
# with new power-flow results ("case with 380 kV done").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: RowNumber, B, C, D, E, F, I, J, K, L, M, N
$data = @(
    @(2, 1.02,1.028639757754075,1.032432101668982,1.032273333469335,1.037794868816087,1.031901358382663,1.033790957529502,1.03523734969345,1.035079039519726,1.040584735341346,1.015108602799619),
    @(3, 1.02,1.030079525528191,1.033516232239543,1.033661736548206,1.039352271275644,1.032222900516538,1.034868909154358,1.036129514706875,1.036274629989866,1.041950042621673,1.015479869650347),
    @(4, 1.02,1.031009655082153,1.034216236757715,1.034559022327542,1.040358703648867,1.032428923691462,1.03556454850967,1.036704720795414,1.037046636580516,1.042831695376559,1.01571911194643),
    @(5, 1.02,1.03140033058385,1.034510163947395,1.03493598330738,1.040781501908025,1.03251505021882,1.035856553516473,1.036946042691687,1.037370805920386,1.043201920249526,1.015819453550102),
    @(6, 1.02,1.031465906350229,1.034559494945747,1.034999261770043,1.040852473888543,1.032529482817859,1.035905556688922,1.036986532810878,1.037425213129551,1.043264058098595,1.015836287564445),
    @(7, 1.02,1.031014876677347,1.034220165615423,1.034564060304422,1.040364354292045,1.032430076424224,1.035568452024223,1.03670794728917,1.037050969641974,1.042836643988559,1.015720453640545),
    @(8, 1.02,1.02912664827153,1.032798802069823,1.032742783155061,1.038321476571445,1.032010447741993,1.03415564618365,1.035539295178763,1.035483433867026,1.041046522761324,1.015234280618789),
    @(9, 1.02,1.025787559417351,1.03028248325455,1.029524722386425,1.034711264242695,1.031255331879469,1.031651561515681,1.033463825811856,1.032708569056878,1.037878051781254,1.014369907225286),
    @(10, 1.02,1.023553113072921,1.028596788016358,1.027373074896923,1.032296959189953,1.03074126449673,1.029972060581969,1.032069047612976,1.030849782298703,1.035755840675975,1.013788391730817),
    @(11, 1.02,1.022583480669071,1.027864865857995,1.026439810279418,1.031249647534746,1.030516113555476,1.029242344650286,1.031462392046729,1.030042720679324,1.03483445123585,1.013535315835277),
    @(12, 1.02,1.022222991191345,1.027592690487887,1.026092908015191,1.030860333861036,1.0304320960611,1.028970916355662,1.031236641063654,1.029742605269121,1.034491828082976,1.013441118374499),
    @(13, 1.02,1.022300332252041,1.027651087007677,1.026167330996788,1.030943856495388,1.030450135622577,1.029029155916811,1.031285084150205,1.029806996397885,1.034565339163727,1.013461332853983),
    @(14, 1.02,1.022553689174958,1.02784237405044,1.026411140298408,1.031217472829939,1.030509176539464,1.029219916078748,1.031443739837694,1.030017919957385,1.034806137663224,1.0135275334097),
    @(15, 1.02,1.02270974742433,1.027960191474315,1.026561326393502,1.031386017434061,1.030545502323392,1.029337399250953,1.031541438076041,1.030147832083448,1.034954451201087,1.013568296012073),
    @(16, 1.02,1.023617419323205,1.028645320554275,1.027434978494671,1.032366424915788,1.030756152963792,1.030020436594655,1.032109251848557,1.030903297479295,1.035816937721347,1.013805160483092),
    @(17, 1.02,1.024186208744636,1.029074543148314,1.027982566993172,1.032980892614796,1.030887602483723,1.030448219032022,1.032464697824432,1.031376588414415,1.036357288455098,1.013953396234237),
    @(18, 1.02,1.024517771750715,1.029324708358431,1.028301813502427,1.033339118413481,1.0309640282717,1.030697498362484,1.032671762587065,1.031652439734687,1.036672229330616,1.014039736683745),
    @(19, 1.02,1.024630792134009,1.029409975652961,1.028410642663397,1.033461233439878,1.030990045751322,1.030782455843699,1.032742322233983,1.031746462282943,1.03677957616965,1.014069155754093),
    @(20, 1.02,1.024125204006333,1.029028511637005,1.027923831769975,1.032914985026232,1.030873524707192,1.030402346786041,1.032426588855945,1.031325830680588,1.036299338468015,1.013937504674358),
    @(21, 1.02,1.0224790909152,1.027786053301504,1.026339351395983,1.031136907875231,1.030491801146736,1.029163752489423,1.031397031111046,1.029955817580324,1.034735239015018,1.013508044363164),
    @(22, 1.02,1.021442230420126,1.027003092009473,1.025341697169245,1.030017245657194,1.030249559568009,1.028382801246211,1.03074732036393,1.029092483862089,1.033749635688362,1.013236903258188),
    @(23, 1.02,1.021992071367079,1.027418325082221,1.025870710577204,1.0306109655918,1.030378189225451,1.028797008815335,1.03109197237124,1.029550340810925,1.03427233338431,1.013380747388891),
    @(24, 1.02,1.024152770046772,1.029049311892383,1.027950372160877,1.032944766393634,1.030879886610178,1.030423075217923,1.032443809464664,1.031348766568707,1.036325524303914,1.013944685766245),
    @(25, 1.02,1.026652238225113,1.030934427200836,1.030357744573867,1.035645874864948,1.031452417133857,1.032300684853697,1.034002326520043,1.033427476121366,1.038698888598835,1.014594288848339),
)

$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

foreach ($row in $data) {
    $rowNum = $row[0]
    for ($i = 0; $i -lt $colIndexes.Count; $i++) {
        $colNum = $colIndexes[$i]
        $value = $row[$i + 1]
        $ws.Cells.Item($rowNum, $colNum).Value = $value
    }
}
